$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# chapter 4 ... ing - add new vocabulary entries (multiplicity / order / zero)
$ws.Range("A95").Value = "multiplicity"
$ws.Range("B95").Value = "중복도"
$ws.Range("D95").Value = "용어사전"

$ws.Range("A96").Value = "order"
$ws.Range("B96").Value = "차수"
$ws.Range("D96").Value = "근의 차수"

$ws.Range("A97").Value = "zero"
$ws.Range("B97").Value = "근"
$ws.Range("D97").Value = "근으로 통일 (해보다는…)"

# Update the view state to match the target (scrolled down, new selection on A98)
$excel.ActiveWindow.ScrollRow = 67
$ws.Range("A98").Select()
